$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = -0.10122516985340724
$ws.Range("B1").Value = 0.10112247416962816
$ws.Range("A2").Value = -0.055012127821452594
$ws.Range("B2").Value = 0.054754468552287605
$ws.Range("A3").Value = 0.048178707084115757
$ws.Range("B3").Value = -0.04824299100872409
$ws.Range("A4").Value = -0.23575105850395772
$ws.Range("B4").Value = 0.23441798158646421
$ws.Range("A5").Value = -0.22841798208705821
$ws.Range("B5").Value = 0.22569808510708356
$ws.Range("A6").Value = -0.10938910939672919
$ws.Range("B6").Value = 0.10920959081280746
$ws.Range("A7").Value = -0.089209591410506661
$ws.Range("B7").Value = 0.088756962029750142
$ws.Range("A8").Value = -0.06875696263328912
$ws.Range("B8").Value = 0.068365534608057388
$ws.Range("A9").Value = -0.06236553514481713
$ws.Range("B9").Value = 0.062032231937378057
$ws.Range("A10").Value = -0.056032232481882716
$ws.Range("B10").Value = 0.055984743502428103
$ws.Range("A11").Value = -0.05148474403984693
$ws.Range("B11").Value = 0.051404383320321045
$ws.Range("A12").Value = -0.045404383867382325
$ws.Range("B12").Value = 0.045154728350481399
$ws.Range("A13").Value = -0.039154728905495872
$ws.Range("B13").Value = 0.039086786850258015
$ws.Range("A14").Value = -0.027086787439746907
$ws.Range("B14").Value = 0.027053985319748897
$ws.Range("A15").Value = -0.021053985878690895
$ws.Range("B15").Value = 0.021028171092475567
$ws.Range("A16").Value = -0.015028171652957667
$ws.Range("B16").Value = 0.015004528111584481
$ws.Range("A17").Value = -0.0090045286741560204
$ws.Range("B17").Value = 0.008999999420990612
$ws.Range("A18").Value = -0.036110346859878462
$ws.Range("B18").Value = 0.036096577391980134
$ws.Range("A19").Value = -0.027096577898665153
$ws.Range("B19").Value = 0.02701368664391568
$ws.Range("A20").Value = -0.018013687154645908
$ws.Range("B20").Value = 0.018004279716734573
$ws.Range("A21").Value = -0.0090042802279866052
$ws.Range("B21").Value = 0.0089999994884486512
$ws.Range("A22").Value = -0.093933176145352348
$ws.Range("B22").Value = 0.093624893454542502
$ws.Range("A23").Value = -0.084624893967234271
$ws.Range("B23").Value = 0.084124893653477528
$ws.Range("A24").Value = -0.042124894350685338
$ws.Range("B24").Value = 0.041999999299422441
$ws.Range("A25").Value = -0.077798689490116146
$ws.Range("B25").Value = 0.07772081778135842
$ws.Range("A26").Value = -0.071720818302999589
$ws.Range("B26").Value = 0.071626672117883317
$ws.Range("A27").Value = -0.065626672641213801
$ws.Range("B27").Value = 0.065326030645288746
$ws.Range("A28").Value = -0.059326031176412997
$ws.Range("B28").Value = 0.059138700843935688
$ws.Range("A29").Value = -0.047138701412329453
$ws.Range("B29").Value = 0.047064582734735794
$ws.Range("A30").Value = -0.037855992536477423
$ws.Range("B30").Value = 0.037762205306250074
$ws.Range("A31").Value = -0.027019079977540272
$ws.Range("B31").Value = 0.027000796273879857
$ws.Range("A32").Value = -0.0060007968972168868
$ws.Range("B32").Value = 0.0059999994570194559
